$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each cell holds a plain text value (inline string) in the source file, e.g. "282.01"
# or "6.23%". Setting .Value directly would let Excel auto-convert these look-alike
# numeric/percent strings into real numbers, so we force a Text number format first,
# write the literal string, then drop back to the Normal style so no stray
# number-format style lingers on the cell (matches the original unstyled text cells).
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '281.71'
Set-TextValue 'E2' '6.00%'
Set-TextValue 'D3' '26.91'
Set-TextValue 'E3' '0.75%'
Set-TextValue 'D4' '4.947'
Set-TextValue 'E4' '5.10%'
Set-TextValue 'E5' '5.52%'
Set-TextValue 'D6' '7.014'
Set-TextValue 'E6' '4.80%'
Set-TextValue 'D7' '3.348'
Set-TextValue 'E7' '5.50%'
Set-TextValue 'D8' '0.8877'
Set-TextValue 'E8' '4.40%'
Set-TextValue 'D9' '1.018'
Set-TextValue 'E9' '12.38%'
Set-TextValue 'D10' '0.1503'
Set-TextValue 'E10' '6.88%'
Set-TextValue 'D11' '0.05232'
Set-TextValue 'E11' '3.29%'
Set-TextValue 'D12' '0.07388'
Set-TextValue 'E12' '3.99%'
Set-TextValue 'D13' '0.03111'
Set-TextValue 'E13' '-1.68%'
Set-TextValue 'D14' '0.09064'
Set-TextValue 'E14' '0.48%'
Set-TextValue 'D15' '0.001557'
Set-TextValue 'E15' '1.12%'
Set-TextValue 'D16' '0.0006321'
Set-TextValue 'E16' '4.44%'
Set-TextValue 'D17' '0.006054'
Set-TextValue 'E17' '1.50%'
Set-TextValue 'D18' '3.499'
Set-TextValue 'E18' '1.23%'
Set-TextValue 'D19' '2.298'
Set-TextValue 'E19' '0.90%'
Set-TextValue 'D20' '0.3115'
Set-TextValue 'E20' '0.90%'
Set-TextValue 'D21' '0.1332'
Set-TextValue 'E21' '2.49%'
Set-TextValue 'D22' '3.921'
Set-TextValue 'E22' '-3.86%'
Set-TextValue 'D23' '0.04363'
Set-TextValue 'E23' '2.95%'
Set-TextValue 'E24' '-0.39%'
Set-TextValue 'D25' '0.003698'
Set-TextValue 'E25' '-10.62%'
Set-TextValue 'D26' '0.0001199'
Set-TextValue 'E26' '-0.17%'
Set-TextValue 'D27' '0.0001694'
Set-TextValue 'E27' '0.65%'
Set-TextValue 'D40' '0.04126'
Set-TextValue 'E40' '5.16%'
Set-TextValue 'D41' '0.006643'
Set-TextValue 'E41' '59.21%'
Set-TextValue 'D42' '0.1180'
Set-TextValue 'E42' '5.81%'
Set-TextValue 'E43' '11.67%'
Set-TextValue 'D44' '0.01257'
Set-TextValue 'E44' '8.67%'
Set-TextValue 'D45' '0.00005263'
Set-TextValue 'E45' '3.05%'
Set-TextValue 'E46' '-0.11%'
Set-TextValue 'E47' '1,902.24%'
Set-TextValue 'D48' '0.02249'
Set-TextValue 'E48' '-8.11%'
Set-TextValue 'E49' '-0.11%'
Set-TextValue 'E50' '-0.18%'
